# Updates for removing premium Power Automate license.
# Adds a new locale row (en-ES) and four new label columns
# (OwnerOrMemberLbl, OwnerLbl, MemberLbl, NoPermissionLbl) to the
# colTranslation table on the "Principal" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# --- New locale row (row 3): en-ES ------------------------------------
# Populate in the same left-to-right order Excel would have used so the
# shared-string table grows in the expected order.
$ws.Range("C3").Value = "`nBonjour,"
$ws.Range("C3").WrapText = $true
$ws.Range("D3").Value = "Nouveau drop-in"
$ws.Range("E3").Value = "Visiter une classe virtuelle"
$ws.Range("F3").Value = "Drop-ins programmés"
$ws.Range("G3").Value = "Afficher les détails ou prolonger la durée"
$ws.Range("A3").Value = "en-ES"
$ws.Range("B3").Value = $true

# Match the row height Excel computed for the wrapped two-line cell.
$ws.Rows.Item(3).RowHeight = 28.8

# --- New columns: OwnerOrMemberLbl / OwnerLbl / MemberLbl / NoPermissionLbl
$newCol = $lo.ListColumns.Add()
$ws.Range("CL1").Value = "OwnerOrMemberLbl"
$ws.Range("CL2").Value = "Dropping the user as :"

$newCol = $lo.ListColumns.Add()
$ws.Range("CM1").Value = "OwnerLbl"
$ws.Range("CM2").Value = "Owner"

$newCol = $lo.ListColumns.Add()
$ws.Range("CN1").Value = "MemberLbl"
$ws.Range("CN2").Value = "Member"

$newCol = $lo.ListColumns.Add()
$ws.Range("CO1").Value = "NoPermissionLbl"
$ws.Range("CO2").Value = "Uh oh! Looks like you don't have permissions to use this app. Please contact IT/support."

# Match the active selection Excel left the workbook in.
$ws.Range("CO1").Select()
